$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("optimization_parameters")

# Remove the extra "value" header cells in C1:F1 (row 1 only needs A1:B1)
$ws.Range("C1:F1").ClearContents()

# Insert a new row after the "production_function" row (row 8) for the new
# "L_curve" optimization parameter, pushing estimate_params..Deletion down by one
$ws.Rows.Item(9).Insert()
$ws.Range("A9").Value = "L_curve"
$ws.Range("B9").Value = 0
$ws.Range("B9").NumberFormat = $ws.Range("B2").NumberFormat

# The old "Deletion" row has shifted from row 16 down to row 17; remove it entirely
$ws.Rows.Item(17).Delete()

# Match the saved active-cell selection on this sheet
$ws.Activate()
$ws.Range("A24").Select()
